$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns (they look numeric but are text)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "31.316.67"
$ws.Range("E2").Value = "  +1.92%  "

$ws.Range("D3").Value = "2.001.33"
$ws.Range("E3").Value = "  +4.53%  "

$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.48%  "

$ws.Range("D5").Value = "0.7756"
$ws.Range("E5").Value = "  +38.32%  "

$ws.Range("D6").Value = "255.75"
$ws.Range("E6").Value = "  +2.66%  "

$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D8").Value = "0.3480"
$ws.Range("E8").Value = "  +16.29%  "

$ws.Range("D9").Value = "28.13"
$ws.Range("E9").Value = "  +19.49%  "

$ws.Range("D10").Value = "0.07165"
$ws.Range("E10").Value = "  +8.31%  "

$ws.Range("D11").Value = "0.8481"
$ws.Range("E11").Value = "  +9.39%  "

$ws.Range("D12").Value = "0.08199"
$ws.Range("E12").Value = "  +4.05%  "

$ws.Range("D13").Value = "101.15"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.001.20"
$ws.Range("E14").Value = "  +6.09%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.653"
$ws.Range("E15").Value = "  +6.63%  "

$ws.Range("D16").Value = "15.49"
$ws.Range("E16").Value = "  +15.65%  "

$ws.Range("D17").Value = "272.87"
$ws.Range("E17").Value = "  -4.05%  "

$ws.Range("D18").Value = "31.307.94"
$ws.Range("E18").Value = "  +2.21%  "

$ws.Range("D19").Value = "0.000008263"
$ws.Range("E19").Value = "  +8.60%  "

$ws.Range("D20").Value = "6.014"
$ws.Range("E20").Value = "  +10.99%  "

$ws.Range("D21").Value = "2.263.11"
$ws.Range("E21").Value = "  +6.13%  "

$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("D24").Value = "7.150"
$ws.Range("E24").Value = "  +9.49%  "

$ws.Range("D25").Value = "10.11"
$ws.Range("E25").Value = "  +9.17%  "

$ws.Range("D26").Value = "164.54"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("D27").Value = "0.1420"
$ws.Range("E27").Value = "  +35.65%  "

$ws.Range("D28").Value = "20.02"
$ws.Range("E28").Value = "  +3.44%  "

$ws.Range("D29").Value = "2.403"
$ws.Range("E29").Value = "  +23.70%  "

$ws.Range("D30").Value = "1.606"
$ws.Range("E30").Value = "  +6.13%  "

$ws.Range("D31").Value = "4.658"
$ws.Range("E31").Value = "  +7.82%  "

$ws.Range("D32").Value = "1.369"
$ws.Range("E32").Value = "  +1.50%  "

$ws.Range("D33").Value = "4.486"
$ws.Range("E33").Value = "  +5.54%  "

$ws.Range("D34").Value = "0.05373"
$ws.Range("E34").Value = "  +9.40%  "

$ws.Range("D35").Value = "1.271"
$ws.Range("E35").Value = "  +10.56%  "

$ws.Range("D36").Value = "0.7896"
$ws.Range("E36").Value = "  +11.11%  "

$ws.Range("D37").Value = "2.776"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("D38").Value = "0.9989"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").Value = "0.02012"
$ws.Range("E39").Value = "  +4.57%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.935"
$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "86.57"
$ws.Range("E41").Value = "  +12.97%  "

$ws.Range("D42").Value = "6.848"
$ws.Range("E42").Value = "  +7.48%  "

$ws.Range("E43").Value = "  +8.02%  "

$ws.Range("D44").Value = "0.4705"
$ws.Range("E44").Value = "  +8.11%  "

$ws.Range("D45").Value = "0.8606"
$ws.Range("E45").Value = "  +1.68%  "

$ws.Range("D46").Value = "105.71"
$ws.Range("E46").Value = "  +4.18%  "

$ws.Range("D47").Value = "10.26"
$ws.Range("E47").Value = "  +1.81%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.0000"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.796"
$ws.Range("E49").Value = "  +8.83%  "

$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "3.107"
$ws.Range("E50").Value = "  +47.25%  "

$ws.Range("D51").Value = "37.90"
$ws.Range("E51").Value = "  +6.74%  "

# Restore default style so no stray number-format style id is left on the cells
$ws.Range("D2:E51").Style = "Normal"

